# Add a new quiz sheet "12_" (before "MultC"), update the old E_gen quiz
# sheet "11_" to use the corrected/shifted comment text, and leave the new
# sheet "12_" as the active/selected tab - matching the upstream commit
# "Add files via upload".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "12_" sheet by copying "11_" (same column widths,
#    same cell styles/fills) and placing it immediately before "MultC".
# ---------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("11_")
$multC    = $wb.Worksheets.Item("MultC")
$srcSheet.Copy($multC)
$newSheet = $wb.ActiveSheet
$newSheet.Name = "12_"

# ---------------------------------------------------------------------
# 2. Populate the new sheet's content (header + 4 answer rows + blank row).
# ---------------------------------------------------------------------
$newSheet.Range("A1").Value = "Which of the following would mean that we *could not* make the simplification 'E_dot = mC dT/dt'?"
$newSheet.Range("B1").Value = "Correct"
$newSheet.Range("C1").Value = "Comment"

$newSheet.Range("A2").Value = "There was phase change in the system"
$newSheet.Range("B2").Value = "Y"
$newSheet.Range("C2").Value = "In this case, the thermal energy of the system would not be fully explained by the change in temperature."

$newSheet.Range("A3").Value = "Mass crossed the boundary of the control volume"
$newSheet.Range("B3").Value = "Y"
$newSheet.Range("C3").Value = "In this case, the energy of the system would change in a way not fully explained by the temperature of the system."

$newSheet.Range("A4").Value = "The control volume accelerated in space, but did not change shape or mass"
$newSheet.Range("B4").Value = "Y"
$newSheet.Range("C4").Value = "In this case, the system would gain (or lose) kinetic energy: so this change could not be calculated using just temperature."

$newSheet.Range("A5").Value = "The temperature in the control volume changed in a nonlinear fashion"
$newSheet.Range("B5").Value = "N"
$newSheet.Range("C5").Value = "The simplification is a instantaneous equation (note the derivative), so it can handle nonlinear changes."

# Row 6 stays blank (carried over from the copied "11_" sheet already has
# the right fill style on A6:C6).

# ---------------------------------------------------------------------
# 3. Row heights for the new sheet: 60 / 45 / 45 / 45 / 45.
# ---------------------------------------------------------------------
$newSheet.Rows.Item(1).RowHeight = 60
$newSheet.Rows.Item(2).RowHeight = 45
$newSheet.Rows.Item(3).RowHeight = 45
$newSheet.Rows.Item(4).RowHeight = 45
$newSheet.Rows.Item(5).RowHeight = 45

# ---------------------------------------------------------------------
# 4. Update the old "11_" sheet: the stray comment string picked up a
#    wording fix ("much" -> "must") that now lives at a new shared-string
#    slot, so re-point C2:C5 at the corrected text in the same order.
# ---------------------------------------------------------------------
$oldSheet = $wb.Worksheets.Item("11_")
$oldSheet.Range("C2").Value = "This is not E_gen: latent heat energy is considered thermal energy, and so must be accounted for in a thermal/ mechanical calculation."
$oldSheet.Range("C3").Value = "Yes, nuclear energy can be redefined as E_gen."
$oldSheet.Range("C4").Value = "Yes, chemical energy can be redefined as E_gen."
$oldSheet.Range("C5").Value = "Yes, electrical energy can be redefined as E_gen."

# Leave "11_" no longer the selected tab; move its cursor to G7 (left over
# UI state from editing), then return focus to the new "12_" sheet.
$oldSheet.Activate()
$oldSheet.Range("G7").Select()

$newSheet.Activate()
$newSheet.Range("C15").Select()
